$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log2")
$ws.Range("A1").Value = "TEST"
